$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 280
$ws.Range("F3").Value = 89
$ws.Range("F4").Value = 1201
$ws.Range("F5").Value = 830
$ws.Range("F6").Value = 861
$ws.Range("F7").Value = 1571
$ws.Range("F8").Value = 321
$ws.Range("F9").Value = 1069
$ws.Range("F11").Value = 81
$ws.Range("F12").Value = 209
$ws.Range("F13").Value = 64
$ws.Range("F14").Value = 532
$ws.Range("F15").Value = 80
$ws.Range("F16").Value = 49
$ws.Range("F17").Value = 15
$ws.Range("F20").Value = 591
$ws.Range("F22").Value = 68
$ws.Range("F24").Value = 790
$ws.Range("F25").Value = 266
$ws.Range("F26").Value = 206
$ws.Range("F28").Value = 379

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1047
$ws.Range("F4").Value = 287
$ws.Range("F6").Value = 186

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 270

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 270
$ws.Range("F3").Value = 280
$ws.Range("F4").Value = 1047
$ws.Range("F5").Value = 89
$ws.Range("F6").Value = 1201
$ws.Range("F7").Value = 830
$ws.Range("F8").Value = 861
$ws.Range("F9").Value = 1571
$ws.Range("F10").Value = 321
$ws.Range("F11").Value = 1069
$ws.Range("F13").Value = 81
$ws.Range("F14").Value = 209
$ws.Range("F15").Value = 64
$ws.Range("F16").Value = 532
$ws.Range("F17").Value = 80
$ws.Range("F18").Value = 49
$ws.Range("F20").Value = 15
$ws.Range("F21").Value = 287
$ws.Range("F25").Value = 186
$ws.Range("F26").Value = 186
$ws.Range("F27").Value = 591
$ws.Range("F29").Value = 68
$ws.Range("F31").Value = 790
$ws.Range("F32").Value = 266
$ws.Range("F34").Value = 206
$ws.Range("F40").Value = 379
